$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 138
$ws.Range("I9").Value = 86.28570999999999
$ws.Range("K9").Value = 86.28570999999999
$ws.Range("M9").Value = 82.71429000000001
$ws.Range("H15").Value = 1815.2646
$ws.Range("I15").Value = 1815.2646
$ws.Range("K15").Value = 5445.793799999999
$ws.Range("M15").Value = -5276.793799999999
$ws.Range("H18").Value = 2990.6667
$ws.Range("I18").Value = 2490.5
$ws.Range("K18").Value = 2490.5
$ws.Range("M18").Value = -2206.5
$ws.Range("H32").Value = 2564.4
$ws.Range("J32").Value = 2691.3333
$ws.Range("L32").Value = 2691.3333
$ws.Range("N32").Value = -3343.3333
$ws.Range("H37").Value = 24670.666
$ws.Range("J37").Value = 1999
$ws.Range("L37").Value = 5997
$ws.Range("N37").Value = -6249
$ws.Range("H38").Value = 1058.6923
$ws.Range("I38").Value = 228.4
$ws.Range("J38").Value = 3826.3333
$ws.Range("K38").Value = 685.2
$ws.Range("L38").Value = 11478.9999
$ws.Range("M38").Value = -313.2
$ws.Range("N38").Value = -12222.9999
$ws.Range("H40").Value = 125002680
$ws.Range("J40").Value = 250001980
$ws.Range("L40").Value = 250001980
$ws.Range("N40").Value = -250002330
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("H70").Value = 1220949.6
$ws.Range("I70").Value = 6098261.5
$ws.Range("J70").Value = 1621.625
$ws.Range("K70").Value = 18294784.5
$ws.Range("L70").Value = 4864.875
$ws.Range("M70").Value = -18294514.5
$ws.Range("N70").Value = -5404.875
$ws.Range("H73").Value = 1220949.6
$ws.Range("I73").Value = 6098261.5
$ws.Range("J73").Value = 1621.625
$ws.Range("K73").Value = 18294784.5
$ws.Range("L73").Value = 4864.875
$ws.Range("M73").Value = -18293848.5
$ws.Range("N73").Value = -6736.875
$ws.Range("H74").Value = 7232.35
$ws.Range("J74").Value = 9516
$ws.Range("L74").Value = 9516
$ws.Range("N74").Value = -11388
$ws.Range("H76").Value = 6692.25
$ws.Range("I76").Value = 2789
$ws.Range("K76").Value = 2789
$ws.Range("M76").Value = -2474
$ws.Range("H77").Value = 7232.35
$ws.Range("J77").Value = 9516
$ws.Range("L77").Value = 47580
$ws.Range("N77").Value = -56940
$ws.Range("H79").Value = 6692.25
$ws.Range("I79").Value = 2789
$ws.Range("K79").Value = 2789
$ws.Range("M79").Value = -1697
$ws.Range("H80").Value = 3735824.5
$ws.Range("I80").Value = 2801806.8
$ws.Range("J80").Value = 4330199.5
$ws.Range("K80").Value = 8405420.399999999
$ws.Range("L80").Value = 12990598.5
$ws.Range("M80").Value = -8404422.399999999
$ws.Range("N80").Value = -12992594.5
$ws.Range("H83").Value = 3735824.5
$ws.Range("I83").Value = 2801806.8
$ws.Range("J83").Value = 4330199.5
$ws.Range("K83").Value = 25216261.2
$ws.Range("L83").Value = 38971795.5
$ws.Range("M83").Value = -25211269.2
$ws.Range("N83").Value = -38981779.5
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("H135").Value = 1351.5358
$ws.Range("I135").Value = 1160.8518
$ws.Range("K135").Value = 10447.6662
$ws.Range("M135").Value = -7912.6662
$ws.Range("H137").Value = 1953.7693
$ws.Range("I137").Value = 1491.5834
$ws.Range("J137").Value = 7500
$ws.Range("K137").Value = 4474.7502
$ws.Range("L137").Value = 22500
$ws.Range("M137").Value = -1924.7502
$ws.Range("N137").Value = -27600
$ws.Range("H138").Value = 5175.793
$ws.Range("I138").Value = 3408.84
$ws.Range("J138").Value = 16219.25
$ws.Range("K138").Value = 10226.52
$ws.Range("L138").Value = 48657.75
$ws.Range("M138").Value = -5086.52
$ws.Range("N138").Value = -58937.75
$ws.Range("H141").Value = 14711362
$ws.Range("I141").Value = 17245092
$ws.Range("J141").Value = 15727.4
$ws.Range("K141").Value = 51735276
$ws.Range("L141").Value = 47182.2
$ws.Range("M141").Value = -51730096
$ws.Range("N141").Value = -57542.2
$ws.Range("M64").ClearContents()
$ws.Range("M67").ClearContents()
$ws.Range("N121").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3532.8333
$ws.Range("I2").Value = 3413.5
$ws.Range("K2").Value = 3413.5
$ws.Range("M2").Value = -3300.5
$ws.Range("H32").Value = 2814.5532
$ws.Range("I32").Value = 2806.7954
$ws.Range("J32").Value = 2928.3333
$ws.Range("K32").Value = 2806.7954
$ws.Range("L32").Value = 2928.3333
$ws.Range("M32").Value = -2519.7954
$ws.Range("N32").Value = -3502.3333
$ws.Range("H61").Value = 956188.0600000001
$ws.Range("I61").Value = 2521.2856
$ws.Range("K61").Value = 2521.2856
$ws.Range("M61").Value = -2309.2856
$ws.Range("H63").Value = 3599.5
$ws.Range("I63").Value = 3599.5
$ws.Range("K63").Value = 3599.5
$ws.Range("M63").Value = -2913.5
$ws.Range("H66").Value = 3599.5
$ws.Range("I66").Value = 3599.5
$ws.Range("K66").Value = 17997.5
$ws.Range("M66").Value = -14565.5
$ws.Range("H74").Value = 1333.591
$ws.Range("J74").Value = 2686
$ws.Range("L74").Value = 2686
$ws.Range("N74").Value = -4434
$ws.Range("H77").Value = 1333.591
$ws.Range("J77").Value = 2686
$ws.Range("L77").Value = 13430
$ws.Range("N77").Value = -22166
$ws.Range("H102").Value = 1739.8948
$ws.Range("I102").Value = 1153.2142
$ws.Range("K102").Value = 1153.2142
$ws.Range("M102").Value = 468.7858000000001
$ws.Range("H103").Value = 120000
$ws.Range("J103").Value = 120000
$ws.Range("L103").Value = 120000
$ws.Range("N103").Value = -122344
$ws.Range("H116").Value = 3532.8333
$ws.Range("I116").Value = 3413.5
$ws.Range("K116").Value = 3413.5
$ws.Range("M116").Value = -1119.5
$ws.Range("H132").Value = 3228618.2
$ws.Range("I132").Value = 2451.4736
$ws.Range("J132").Value = 8336716
$ws.Range("K132").Value = 7354.4208
$ws.Range("L132").Value = 25010148
$ws.Range("M132").Value = -4824.4208
$ws.Range("N132").Value = -25015208
$ws.Range("H136").Value = 956188.0600000001
$ws.Range("I136").Value = 2521.2856
$ws.Range("K136").Value = 7563.8568
$ws.Range("M136").Value = -5013.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3532.8333
$ws.Range("I3").Value = 3413.5
$ws.Range("K3").Value = 3413.5
$ws.Range("M3").Value = -3299.5
$ws.Range("H86").Value = 3734.4211
$ws.Range("I86").Value = 1228.4445
$ws.Range("J86").Value = 5989.8
$ws.Range("K86").Value = 1228.4445
$ws.Range("L86").Value = 5989.8
$ws.Range("M86").Value = -105.4445000000001
$ws.Range("N86").Value = -8235.799999999999
$ws.Range("H89").Value = 3734.4211
$ws.Range("I89").Value = 1228.4445
$ws.Range("J89").Value = 5989.8
$ws.Range("K89").Value = 6142.2225
$ws.Range("L89").Value = 29949
$ws.Range("M89").Value = -526.2224999999999
$ws.Range("N89").Value = -41181
$ws.Range("H107").Value = 5323.846
$ws.Range("I107").Value = 6119.727
$ws.Range("K107").Value = 6119.727
$ws.Range("M107").Value = -4199.727
$ws.Range("H132").Value = 130053
$ws.Range("J132").Value = 130053
$ws.Range("L132").Value = 130053
$ws.Range("N132").Value = -140173
$ws.Range("H134").Value = 4547838.5
$ws.Range("I134").Value = 2122.1
$ws.Range("K134").Value = 6366.299999999999
$ws.Range("M134").Value = -3831.299999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2843.24
$ws.Range("I31").Value = 3396.3572
$ws.Range("J31").Value = 2139.2727
$ws.Range("K31").Value = 3396.3572
$ws.Range("L31").Value = 2139.2727
$ws.Range("M31").Value = -3101.3572
$ws.Range("N31").Value = -2729.2727
$ws.Range("H34").Value = 2843.24
$ws.Range("I34").Value = 3396.3572
$ws.Range("J34").Value = 2139.2727
$ws.Range("K34").Value = 3396.3572
$ws.Range("L34").Value = 2139.2727
$ws.Range("M34").Value = -3194.3572
$ws.Range("N34").Value = -2543.2727
$ws.Range("H58").Value = 1909.579
$ws.Range("I58").Value = 1303.8182
$ws.Range("J58").Value = 2742.5
$ws.Range("K58").Value = 1303.8182
$ws.Range("L58").Value = 2742.5
$ws.Range("M58").Value = -1100.8182
$ws.Range("N58").Value = -3148.5
$ws.Range("H110").Value = 117675.5
$ws.Range("J110").Value = 117675.5
$ws.Range("L110").Value = 117675.5
$ws.Range("N110").Value = -125855.5
$ws.Range("H122").Value = 4532.923
$ws.Range("I122").Value = 5529.4
$ws.Range("J122").Value = 3910.125
$ws.Range("K122").Value = 16588.2
$ws.Range("L122").Value = 11730.375
$ws.Range("M122").Value = -14138.2
$ws.Range("N122").Value = -16630.375
$ws.Range("H132").Value = 2527.4285
$ws.Range("I132").Value = 2408.5454
$ws.Range("K132").Value = 7225.6362
$ws.Range("M132").Value = -4695.6362
$ws.Range("H136").Value = 1909.579
$ws.Range("I136").Value = 1303.8182
$ws.Range("J136").Value = 2742.5
$ws.Range("K136").Value = 3911.4546
$ws.Range("L136").Value = 8227.5
$ws.Range("M136").Value = -1361.4546
$ws.Range("N136").Value = -13327.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 11114202
$ws.Range("J12").Value = 15155722
$ws.Range("L12").Value = 45467166
$ws.Range("N12").Value = -45467512
$ws.Range("H68").Value = 1265.579
$ws.Range("I68").Value = 1299.2
$ws.Range("J68").Value = 1253.5714
$ws.Range("K68").Value = 3897.6
$ws.Range("L68").Value = 3760.7142
$ws.Range("M68").Value = -3086.6
$ws.Range("N68").Value = -5382.7142
$ws.Range("H71").Value = 1265.579
$ws.Range("I71").Value = 1299.2
$ws.Range("J71").Value = 1253.5714
$ws.Range("K71").Value = 11692.8
$ws.Range("L71").Value = 11282.1426
$ws.Range("M71").Value = -7636.800000000001
$ws.Range("N71").Value = -19394.1426
$ws.Range("H93").Value = 13483
$ws.Range("J93").Value = 17644
$ws.Range("L93").Value = 52932
$ws.Range("N93").Value = -56676
$ws.Range("H107").Value = 2847827.5
$ws.Range("J107").Value = 4554903.5
$ws.Range("L107").Value = 13664710.5
$ws.Range("N107").Value = -13668550.5
$ws.Range("H109").Value = 14708160
$ws.Range("I109").Value = 15625336
$ws.Range("K109").Value = 46876008
$ws.Range("M109").Value = -46874968
$ws.Range("H129").Value = 1321262.2
$ws.Range("I129").Value = 6057.25
$ws.Range("J129").Value = 1671983.6
$ws.Range("K129").Value = 18171.75
$ws.Range("L129").Value = 5015950.800000001
$ws.Range("M129").Value = -13171.75
$ws.Range("N129").Value = -5025950.800000001
$ws.Range("H137").Value = 10805
$ws.Range("J137").Value = 19110
$ws.Range("L137").Value = 57330
$ws.Range("N137").Value = -67530
$ws.Range("H139").Value = 3039.3333
$ws.Range("I139").Value = 3257.5
$ws.Range("J139").Value = 2995.7
$ws.Range("K139").Value = 9772.5
$ws.Range("L139").Value = 8987.099999999999
$ws.Range("M139").Value = -4632.5
$ws.Range("N139").Value = -19267.1
$ws.Range("H141").Value = 6138.5
$ws.Range("J141").Value = 33333
$ws.Range("L141").Value = 99999
$ws.Range("N141").Value = -110359

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 991.1539
$ws.Range("I2").Value = 1088
$ws.Range("J2").Value = 458.5
$ws.Range("K2").Value = 1088
$ws.Range("L2").Value = 458.5
$ws.Range("M2").Value = -975
$ws.Range("N2").Value = -684.5
$ws.Range("H43").Value = 28999
$ws.Range("I43").Value = 28999
$ws.Range("K43").Value = 28999
$ws.Range("M43").Value = -28848
$ws.Range("H57").Value = 27800
$ws.Range("J57").Value = 27800
$ws.Range("L57").Value = 27800
$ws.Range("N57").Value = -29440
$ws.Range("H70").Value = 7819.875
$ws.Range("J70").Value = 7897.778
$ws.Range("L70").Value = 7897.778
$ws.Range("N70").Value = -8437.778
$ws.Range("H73").Value = 7819.875
$ws.Range("J73").Value = 7897.778
$ws.Range("L73").Value = 7897.778
$ws.Range("N73").Value = -9769.778
$ws.Range("H80").Value = 2357.8462
$ws.Range("I80").Value = 1571.9
$ws.Range("K80").Value = 1571.9
$ws.Range("M80").Value = -573.9000000000001
$ws.Range("H82").Value = 74999
$ws.Range("J82").Value = 74999
$ws.Range("L82").Value = 74999
$ws.Range("N82").Value = -75765
$ws.Range("H83").Value = 2357.8462
$ws.Range("I83").Value = 1571.9
$ws.Range("K83").Value = 7859.5
$ws.Range("M83").Value = -2867.5
$ws.Range("H85").Value = 74999
$ws.Range("J85").Value = 74999
$ws.Range("L85").Value = 74999
$ws.Range("N85").Value = -77651
$ws.Range("H92").Value = 58747.2
$ws.Range("J92").Value = 58747.2
$ws.Range("L92").Value = 58747.2
$ws.Range("N92").Value = -62491.2
$ws.Range("H107").Value = 860.05554
$ws.Range("J107").Value = 894.4286
$ws.Range("L107").Value = 894.4286
$ws.Range("N107").Value = -4734.4286
$ws.Range("H132").Value = 4169570.5
$ws.Range("I132").Value = 2800.5625
$ws.Range("J132").Value = 12503110
$ws.Range("K132").Value = 8401.6875
$ws.Range("L132").Value = 37509330
$ws.Range("M132").Value = -5871.6875
$ws.Range("N132").Value = -37514390

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1880
$ws.Range("I7").Value = 1880
$ws.Range("K7").Value = 1880
$ws.Range("M7").Value = -1768
$ws.Range("H68").Value = 4632136
$ws.Range("I68").Value = 6946224
$ws.Range("K68").Value = 6946224
$ws.Range("M68").Value = -6945475
$ws.Range("H71").Value = 4632136
$ws.Range("I71").Value = 6946224
$ws.Range("K71").Value = 34731120
$ws.Range("M71").Value = -34727376
$ws.Range("H126").Value = 1880
$ws.Range("I126").Value = 1880
$ws.Range("K126").Value = 5640
$ws.Range("M126").Value = -3170
$ws.Range("H136").Value = 2493.3438
$ws.Range("I136").Value = 2369.3845
$ws.Range("J136").Value = 3030.5
$ws.Range("K136").Value = 7108.1535
$ws.Range("L136").Value = 9091.5
$ws.Range("M136").Value = -4558.1535
$ws.Range("N136").Value = -14191.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 4282.1816
$ws.Range("I14").Value = 4282.1816
$ws.Range("K14").Value = 4282.1816
$ws.Range("M14").Value = -4114.1816
$ws.Range("H15").Value = 32486.5
$ws.Range("I15").Value = 29998
$ws.Range("K15").Value = 29998
$ws.Range("M15").Value = -29710
$ws.Range("H16").Value = 74999.5
$ws.Range("J16").Value = 74999.5
$ws.Range("L16").Value = 74999.5
$ws.Range("N16").Value = -75583.5
$ws.Range("H19").Value = 14998
$ws.Range("J19").Value = 14998
$ws.Range("L19").Value = 14998
$ws.Range("N19").Value = -15346
$ws.Range("H132").Value = 264770.1
$ws.Range("I132").Value = 1592.9333
$ws.Range("K132").Value = 4778.7999
$ws.Range("M132").Value = -2248.7999
$ws.Range("H136").Value = 195681.8
$ws.Range("I136").Value = 7254.7856
$ws.Range("J136").Value = 915130.4
$ws.Range("K136").Value = 21764.3568
$ws.Range("L136").Value = 2745391.2
$ws.Range("M136").Value = -19214.3568
$ws.Range("N136").Value = -2750491.2
